$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.9432543333333333
$ws.Range("H2").Value = 2.829763
$ws.Range("I2").Value = 0.00840968845628655
$ws.Range("J2").Value = 0.008409688456286549
$ws.Range("M2").Value = 9.084137666666667
$ws.Range("N2").Value = 27.252413
$ws.Range("O2").Value = 0.2765376761551382
$ws.Range("P2").Value = 0.2765376761551382
$ws.Range("Q2").Value = 8.568652218679889
$ws.Range("R2").Value = 77.117869968119
$ws.Range("S2").Value = 0.002325595702890174
$ws.Range("T2").Value = 0.002325595702890174

# Row 3
$ws.Range("G3").Value = 0.9432543333333333
$ws.Range("H3").Value = 2.829763
$ws.Range("I3").Value = 0.00840968845628655
$ws.Range("J3").Value = 0.008409688456286549
$ws.Range("O3").Value = 0.3707916163717078
$ws.Range("P3").Value = 0.3707916163717078
$ws.Range("Q3").Value = 11.48915565671033
$ws.Range("R3").Value = 103.402400910393
$ws.Range("S3").Value = 0.003118241975888982
$ws.Range("T3").Value = 0.003118241975888982

# Row 4
$ws.Range("G4").Value = 0.9432543333333333
$ws.Range("H4").Value = 2.829763
$ws.Range("I4").Value = 0.00840968845628655
$ws.Range("J4").Value = 0.008409688456286549
$ws.Range("M4").Value = 11.58507333333333
$ws.Range("N4").Value = 34.75522
$ws.Range("O4").Value = 0.3526707074731541
$ws.Range("P4").Value = 0.3526707074731541
$ws.Range("Q4").Value = 10.92767062365111
$ws.Range("R4").Value = 98.34903561285999
$ws.Range("S4").Value = 0.002965850777507395
$ws.Range("T4").Value = 0.002965850777507394

# Row 5
$ws.Range("I5").Value = 0.9472012688364344
$ws.Range("J5").Value = 0.9472012688364343
$ws.Range("M5").Value = 9.084137666666667
$ws.Range("N5").Value = 27.252413
$ws.Range("O5").Value = 0.2765376761551382
$ws.Range("P5").Value = 0.2765376761551382
$ws.Range("Q5").Value = 965.1056987354311
$ws.Range("R5").Value = 8685.951288618879
$ws.Range("S5").Value = 0.2619368377352259
$ws.Range("T5").Value = 0.2619368377352259

# Row 6
$ws.Range("I6").Value = 0.9472012688364344
$ws.Range("J6").Value = 0.9472012688364343
$ws.Range("O6").Value = 0.3707916163717078
$ws.Range("P6").Value = 0.3707916163717078
$ws.Range("S6").Value = 0.3512142895011941
$ws.Range("T6").Value = 0.351214289501194

# Row 7
$ws.Range("I7").Value = 0.9472012688364344
$ws.Range("J7").Value = 0.9472012688364343
$ws.Range("M7").Value = 11.58507333333333
$ws.Range("N7").Value = 34.75522
$ws.Range("O7").Value = 0.3526707074731541
$ws.Range("P7").Value = 0.3526707074731541
$ws.Range("Q7").Value = 1230.80700717414
$ws.Range("R7").Value = 11077.26306456726
$ws.Range("S7").Value = 0.3340501416000146
$ws.Range("T7").Value = 0.3340501416000145

# Row 8
$ws.Range("G8").Value = 4.9788
$ws.Range("H8").Value = 14.9364
$ws.Range("I8").Value = 0.04438904270727917
$ws.Range("J8").Value = 0.04438904270727916
$ws.Range("M8").Value = 9.084137666666667
$ws.Range("N8").Value = 27.252413
$ws.Range("O8").Value = 0.2765376761551382
$ws.Range("P8").Value = 0.2765376761551382
$ws.Range("Q8").Value = 45.2281046148
$ws.Range("R8").Value = 407.0529415332
$ws.Range("S8").Value = 0.01227524271702217
$ws.Range("T8").Value = 0.01227524271702216

# Row 9
$ws.Range("G9").Value = 4.9788
$ws.Range("H9").Value = 14.9364
$ws.Range("I9").Value = 0.04438904270727917
$ws.Range("J9").Value = 0.04438904270727916
$ws.Range("O9").Value = 0.3707916163717078
$ws.Range("P9").Value = 0.3707916163717078
$ws.Range("Q9").Value = 60.6434618556
$ws.Range("R9").Value = 545.7911567004
$ws.Range("S9").Value = 0.01645908489462481
$ws.Range("T9").Value = 0.01645908489462481

# Row 10
$ws.Range("G10").Value = 4.9788
$ws.Range("H10").Value = 14.9364
$ws.Range("I10").Value = 0.04438904270727917
$ws.Range("J10").Value = 0.04438904270727916
$ws.Range("M10").Value = 11.58507333333333
$ws.Range("N10").Value = 34.75522
$ws.Range("O10").Value = 0.3526707074731541
$ws.Range("P10").Value = 0.3526707074731541
$ws.Range("Q10").Value = 57.679763112
$ws.Range("R10").Value = 519.117868008
$ws.Range("S10").Value = 0.0156547150956322
$ws.Range("T10").Value = 0.01565471509563219

Write-Output "Updated Igf1-Igf1r values"